$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-02-04 Tuesday" "2025-02-05 Wednesday"

Replace-Text "126×7=" "633×4="
Replace-Text "889×7=" "476×5="
Replace-Text "360×9=" "640×8="
Replace-Text "602×7=" "887×7="
Replace-Text "969×4=" "266×9="

Replace-Text "567×5=" "846×7="
Replace-Text "685×5=" "305×9="
Replace-Text "177×9=" "781×3="
Replace-Text "406×3=" "666×6="
Replace-Text "702×9=" "490×9="

Replace-Text "531×4=" "880×3="
Replace-Text "991×2=" "552×7="
Replace-Text "241×4=" "193×2="
Replace-Text "735×3=" "499×3="
Replace-Text "585×5=" "368×3="

Replace-Text "357×9=" "346×6="
Replace-Text "451×7=" "139×4="
Replace-Text "832×9=" "233×7="
Replace-Text "525×6=" "690×8="
Replace-Text "583×3=" "816×7="

Replace-Text "338×2=" "788×5="
Replace-Text "151×9=" "972×7="
Replace-Text "176×5=" "321×6="
Replace-Text "912×4=" "364×8="
Replace-Text "483×2=" "384×3="
